$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.8794496666666666
$ws.Range("N2").Value = 2.638349
$ws.Range("O2").Value = 0.3488427963707166
$ws.Range("P2").Value = 0.3488427963707166
$ws.Range("Q2").Value = 0.4200368867955556
$ws.Range("R2").Value = 3.78033198116
$ws.Range("S2").Value = 0.3488427963707166
$ws.Range("T2").Value = 0.3488427963707166

# Row 3 updates
$ws.Range("O3").Value = 0.2822103394539786
$ws.Range("P3").Value = 0.2822103394539786
$ws.Range("S3").Value = 0.2822103394539786
$ws.Range("T3").Value = 0.2822103394539786

# Row 4 updates
$ws.Range("M4").Value = 0.8404543333333333
$ws.Range("N4").Value = 2.521363
$ws.Range("O4").Value = 0.3333748945214069
$ws.Range("P4").Value = 0.3333748945214068
$ws.Range("Q4").Value = 0.4014121956577778
$ws.Range("R4").Value = 3.612709760920001
$ws.Range("S4").Value = 0.3333748945214069
$ws.Range("T4").Value = 0.3333748945214068

# Row 5 updates
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08967866666666667
$ws.Range("N5").Value = 0.269036
$ws.Range("O5").Value = 0.035571969653898
$ws.Range("P5").Value = 0.03557196965389799
$ws.Range("Q5").Value = 0.04283172691555556
$ws.Range("R5").Value = 0.38548554224
$ws.Range("S5").Value = 0.035571969653898
$ws.Range("T5").Value = 0.03557196965389799
